$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '26.746.44'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +1.47%  '
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.720.02'
$c.Style = "Normal"
$ws.Range("E3").Value = '  +0.26%  '
$ws.Range("E4").Value = '  +0.33%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '239.59'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.57%  '
$ws.Range("E6").Value = '  +0.39%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.4753'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -2.33%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.2548'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -1.40%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.06104'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -0.97%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '1.719.45'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -0.15%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '15.81'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +2.21%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.06888'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -1.00%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.5941'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -0.43%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '4.394'
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '76.24'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -0.27%  '
$ws.Range("E16").Value = '  +0.41%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '26.589.38'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +0.89%  '
$ws.Range("E18").Value = '  +0.33%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.000007010'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -1.08%  '
$ws.Range("E20").Value = '  -0.11%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '1.941.54'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -0.07%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '4.355'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -1.11%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '8.294'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -1.74%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '5.022'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.53%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '140.23'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +2.68%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '15.05'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -0.98%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '1.775'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +2.67%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '1.382'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -2.09%  '
$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '105.68'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -0.03%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '3.923'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +1.32%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.07853'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -1.23%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '3.604'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -0.11%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.04522'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +2.20%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '2.590'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -0.76%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.9876'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -0.68%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.6085'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -1.41%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.9167'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -1.77%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '2.483'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +4.67%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '1.949'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -1.19%  '
$ws.Range("E40").Value = '  +0.42%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '5.683'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +4.55%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.01475'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +0.24%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '100.20'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +1.06%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.3769'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -0.91%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '6.684'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -2.01%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.1139'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -0.93%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.05342'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -0.15%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '7.764'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +0.36%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '29.55'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -2.70%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '1.224'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +0.84%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '1.003'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +0.31%  '
